$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma-separated names to use periods ---
$ws.Range("E87").Value = "FERNANDEZ. MARIO HUGO"
$ws.Range("E164").Value = "FERNANDEZ. MARIO HUGO"
$ws.Range("E89").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E169").Value = "DODERA. JORGE ABELARDO"
$ws.Range("E198").Value = "DODERA. JORGE ABELARDO"
$ws.Range("E174").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("F135").Value = "MERCANZINI. GASTON ARIEL"

# --- Fix number formatting in column H (Importe): Latin format -> plain decimal ---
$rng = $ws.Range("H2:H223")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$arr = New-Object 'object[,]' 222,1
$arr[0,0] = "20205.00"
$arr[1,0] = "2550.00"
$arr[2,0] = "48200.00"
$arr[3,0] = "1707.00"
$arr[4,0] = "54738.98"
$arr[5,0] = "2400.00"
$arr[6,0] = "4580.00"
$arr[7,0] = "90.00"
$arr[8,0] = "60450.00"
$arr[9,0] = "438219.76"
$arr[10,0] = "12960.00"
$arr[11,0] = "38495.56"
$arr[12,0] = "8902.00"
$arr[13,0] = "410.00"
$arr[14,0] = "4668.70"
$arr[15,0] = "52057.40"
$arr[16,0] = "1772.00"
$arr[17,0] = "10748.26"
$arr[18,0] = "29427.58"
$arr[19,0] = "2700.00"
$arr[20,0] = "5000.00"
$arr[21,0] = "380.00"
$arr[22,0] = "1450.00"
$arr[23,0] = "7.72"
$arr[24,0] = "2085.80"
$arr[25,0] = "14305.50"
$arr[26,0] = "1700.00"
$arr[27,0] = "86720.00"
$arr[28,0] = "9705.00"
$arr[29,0] = "9224.72"
$arr[30,0] = "2333.19"
$arr[31,0] = "3232.02"
$arr[32,0] = "47568.00"
$arr[33,0] = "45044.00"
$arr[34,0] = "70.00"
$arr[35,0] = "2618.00"
$arr[36,0] = "469.25"
$arr[37,0] = "3915.95"
$arr[38,0] = "3371.81"
$arr[39,0] = "4899.00"
$arr[40,0] = "7380.00"
$arr[41,0] = "9877.00"
$arr[42,0] = "817.24"
$arr[43,0] = "13675.00"
$arr[44,0] = "130156.64"
$arr[45,0] = "162.00"
$arr[46,0] = "473.84"
$arr[47,0] = "275.00"
$arr[48,0] = "550.00"
$arr[49,0] = "4108.40"
$arr[50,0] = "39569.04"
$arr[51,0] = "68.50"
$arr[52,0] = "3344.41"
$arr[53,0] = "38.50"
$arr[54,0] = "629.37"
$arr[55,0] = "380.00"
$arr[56,0] = "117679.76"
$arr[57,0] = "256776.62"
$arr[58,0] = "9590.00"
$arr[59,0] = "12092.41"
$arr[60,0] = "776.55"
$arr[61,0] = "6146.00"
$arr[62,0] = "2629.27"
$arr[63,0] = "7536.00"
$arr[64,0] = "830.00"
$arr[65,0] = "9000.00"
$arr[66,0] = "110.00"
$arr[67,0] = "986.79"
$arr[68,0] = "14181.68"
$arr[69,0] = "3840.00"
$arr[70,0] = "10450.76"
$arr[71,0] = "36822.00"
$arr[72,0] = "450.00"
$arr[73,0] = "30042.50"
$arr[74,0] = "6064.00"
$arr[75,0] = "67330.00"
$arr[76,0] = "3450.30"
$arr[77,0] = "5530.00"
$arr[78,0] = "3360.83"
$arr[79,0] = "2315.00"
$arr[80,0] = "14700.00"
$arr[81,0] = "11100.00"
$arr[82,0] = "950.00"
$arr[83,0] = "567.11"
$arr[84,0] = "3290.00"
$arr[85,0] = "340.00"
$arr[86,0] = "28740.00"
$arr[87,0] = "14560.00"
$arr[88,0] = "1093.95"
$arr[89,0] = "200.00"
$arr[90,0] = "4000.00"
$arr[91,0] = "875.00"
$arr[92,0] = "7722.00"
$arr[93,0] = "371897.61"
$arr[94,0] = "53261.49"
$arr[95,0] = "100.00"
$arr[96,0] = "30.78"
$arr[97,0] = "3575.00"
$arr[98,0] = "1383.28"
$arr[99,0] = "508.64"
$arr[100,0] = "28268.41"
$arr[101,0] = "24.50"
$arr[102,0] = "288.00"
$arr[103,0] = "8657.50"
$arr[104,0] = "54.45"
$arr[105,0] = "1974.00"
$arr[106,0] = "2000.00"
$arr[107,0] = "165.00"
$arr[108,0] = "1840.50"
$arr[109,0] = "577.50"
$arr[110,0] = "533.00"
$arr[111,0] = "3646.00"
$arr[112,0] = "126.00"
$arr[113,0] = "2740.00"
$arr[114,0] = "6200.00"
$arr[115,0] = "7609.30"
$arr[116,0] = "4141.00"
$arr[117,0] = "1338.02"
$arr[118,0] = "49.43"
$arr[119,0] = "5098.50"
$arr[120,0] = "119.00"
$arr[121,0] = "104.50"
$arr[122,0] = "4550.00"
$arr[123,0] = "460.70"
$arr[124,0] = "3300.00"
$arr[125,0] = "40.00"
$arr[126,0] = "1944.00"
$arr[127,0] = "4212.00"
$arr[128,0] = "3818.00"
$arr[129,0] = "6430.00"
$arr[130,0] = "800.00"
$arr[131,0] = "9300.00"
$arr[132,0] = "2685.00"
$arr[133,0] = "18000.00"
$arr[134,0] = "3400.00"
$arr[135,0] = "7950.00"
$arr[136,0] = "26161.23"
$arr[137,0] = "4640.00"
$arr[138,0] = "6582.27"
$arr[139,0] = "8212.00"
$arr[140,0] = "926.40"
$arr[141,0] = "4949.22"
$arr[142,0] = "3075.00"
$arr[143,0] = "12718.17"
$arr[144,0] = "205025.00"
$arr[145,0] = "9360.92"
$arr[146,0] = "2300.00"
$arr[147,0] = "800.00"
$arr[148,0] = "1500.00"
$arr[149,0] = "30492.00"
$arr[150,0] = "3146.00"
$arr[151,0] = "1800.00"
$arr[152,0] = "1105.00"
$arr[153,0] = "1440.00"
$arr[154,0] = "800.00"
$arr[155,0] = "1000.00"
$arr[156,0] = "1568.25"
$arr[157,0] = "1200.00"
$arr[158,0] = "5656.40"
$arr[159,0] = "600.00"
$arr[160,0] = "800.00"
$arr[161,0] = "9280.00"
$arr[162,0] = "7350.00"
$arr[163,0] = "1250.00"
$arr[164,0] = "2299.00"
$arr[165,0] = "1800.00"
$arr[166,0] = "2541.00"
$arr[167,0] = "550.00"
$arr[168,0] = "950.00"
$arr[169,0] = "945.01"
$arr[170,0] = "1100.00"
$arr[171,0] = "237.18"
$arr[172,0] = "3632.00"
$arr[173,0] = "9069.84"
$arr[174,0] = "1193.00"
$arr[175,0] = "5100.00"
$arr[176,0] = "187.00"
$arr[177,0] = "536.00"
$arr[178,0] = "17300.00"
$arr[179,0] = "1847.63"
$arr[180,0] = "2118.31"
$arr[181,0] = "3481.00"
$arr[182,0] = "5126.00"
$arr[183,0] = "3869.00"
$arr[184,0] = "270.00"
$arr[185,0] = "16073.50"
$arr[186,0] = "721.72"
$arr[187,0] = "1223.79"
$arr[188,0] = "3600.00"
$arr[189,0] = "10680.00"
$arr[190,0] = "1631.25"
$arr[191,0] = "4520.00"
$arr[192,0] = "25472.39"
$arr[193,0] = "3500.00"
$arr[194,0] = "2035.51"
$arr[195,0] = "84200.00"
$arr[196,0] = "600.00"
$arr[197,0] = "894489.72"
$arr[198,0] = "12800.00"
$arr[199,0] = "165800.21"
$arr[200,0] = "199334.95"
$arr[201,0] = "238500.00"
$arr[202,0] = "395386.00"
$arr[203,0] = "32500.00"
$arr[204,0] = "63550.00"
$arr[205,0] = "120370.25"
$arr[206,0] = "350141.20"
$arr[207,0] = "303315.00"
$arr[208,0] = "65000.00"
$arr[209,0] = "215241.00"
$arr[210,0] = "217700.00"
$arr[211,0] = "92780.00"
$arr[212,0] = "100000.00"
$arr[213,0] = "3999.99"
$arr[214,0] = "8587.25"
$arr[215,0] = "37110.00"
$arr[216,0] = "17300.00"
$arr[217,0] = "1100.00"
$arr[218,0] = "13500.00"
$arr[219,0] = "248000.00"
$arr[220,0] = "12540.00"
$arr[221,0] = "4900.00"
$rng.Value = $arr
$rng.Style = $origStyle
